$wb  = $excel.ActiveWorkbook
$wsExt = $wb.Worksheets.Item("ExtTest10mm")
$wsFlx = $wb.Worksheets.Item("FlxTest10mm")

# ---------------------------------------------------------------------
# ExtTest10mm (sheet1) - add "actual" / measured routing-path columns
# ---------------------------------------------------------------------
$wsExt.Range("D3").Value = "actual"
$wsExt.Range("D3").Font.Bold = $true
$wsExt.Range("D3").HorizontalAlignment = -4152
$wsExt.Range("E3").Value = 440

$wsExt.Range("D4").Value = "actual"
$wsExt.Range("D4").Font.Bold = $true
$wsExt.Range("D4").HorizontalAlignment = -4152
$wsExt.Range("E4").Value = 55

# Rename the Vas_Pam insertion reference point and update the measured
# coordinates beneath it
$wsExt.Range("J21").Value = "Vas_Pam insertion wrt proximal ring"
$wsExt.Range("J23").Value = 0.02053
$wsExt.Range("K23").Value = -0.07557

# Remove the now-duplicated second block (rows 25-27)
$wsExt.Range("J25:L27").ClearContents() | Out-Null

# ---------------------------------------------------------------------
# FlxTest10mm (sheet2) - same "actual" column plus a new Tendon Length row
# ---------------------------------------------------------------------
$wsFlx.Range("D3").Value = "actual"
$wsFlx.Range("D3").Font.Bold = $true
$wsFlx.Range("D3").HorizontalAlignment = -4152
$wsFlx.Range("E3").Value = 350

$wsFlx.Range("B4").Value = "Tendon Length"
$wsFlx.Range("B4").Font.Bold = $true
$wsFlx.Range("C4").Value = 11
$wsFlx.Range("D4").Value = "?"

# ---------------------------------------------------------------------
# Chart axis rescale on the FlxTest10mm chart (Knee angle vs Torque)
# ---------------------------------------------------------------------
$chart = $wsFlx.ChartObjects(1).Chart
$axCat = $chart.Axes(1)
$axVal = $chart.Axes(2)
$axCat.MinimumScale = -120
$axCat.MaximumScale = 20
$axVal.MinimumScale = -16
$axVal.MaximumScale = 0

# ---------------------------------------------------------------------
# Switch the active sheet / selections to match the saved view state
# ---------------------------------------------------------------------
$wsExt.Range("J24").Select() | Out-Null
$wsFlx.Activate() | Out-Null
$wsFlx.Range("C13:L13").Select() | Out-Null
